$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the data (row 19),
# pushing the previously-existing rows 19-24 down to rows 20-25.
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the latest week's figures.
$ws.Range("A19").Value2 = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value2 = 44510
$ws.Range("E19").Value2 = 8
$ws.Range("F19").Value2 = 100114007
$ws.Range("G19").Value = "Jengibre"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value2 = 40
$ws.Range("K19").Value2 = 15000
$ws.Range("L19").Value2 = 16000
$ws.Range("M19").Value2 = 15500
$ws.Range("N19").Value = "$/caja 13 kilos"
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value2 = 1192
$ws.Range("Q19").Value2 = 13
$ws.Range("R19").Value = "Hortaliza"
